$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "OTROS" column (K), which shifts
# OTROS / A PAGAR / DESCUENTO / PAGOS / SALDO FINAL one column to the right.
$ws.Columns.Item(11).Insert()

# Copy the formatting from the neighboring "JURIDICO" column (J) onto the
# newly inserted column so the new header/data cells pick up the same
# styles already present in the workbook (border, fill, font, number
# format) instead of creating brand-new style entries.
$ws.Range("J1:J6").Copy()
$ws.Range("K1:K6").PasteSpecial(-4122)

# New column header and balance values.
$ws.Cells.Item(1, 11).Value = "SALDO ANTERIOR"
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 11).Value = 0.0
}

# Width for the new "SALDO ANTERIOR" column.
$ws.Columns.Item(11).ColumnWidth = 18.63
